# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns for rows 2-51.
# Values that look like plain numbers are written with a leading apostrophe so
# they stay text cells (matching the original inline-string cells) instead of
# being auto-coerced to numeric values by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.926.97'
$ws.Range("E2").Value = '  -1.15%  '
$ws.Range("D3").Value = '2.043.44'
$ws.Range("E3").Value = '  -2.20%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '''251.42'
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("D6").Value = '''0.668'
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("D7").Value = '''58.76'
$ws.Range("E7").Value = '  +7.71%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '''60.89'
$ws.Range("E9").Value = '  -2.50%  '
$ws.Range("E10").Value = '  +0.74%  '
$ws.Range("D11").Value = '''0.0787'
$ws.Range("E11").Value = '  +4.34%  '
$ws.Range("E12").Value = '  +1.88%  '
$ws.Range("D13").Value = '''16.17'
$ws.Range("E13").Value = '  +4.31%  '
$ws.Range("D14").Value = '2.339.91'
$ws.Range("E14").Value = '  -2.25%  '
$ws.Range("E15").Value = '  -7.09%  '
$ws.Range("E16").Value = '  +6.58%  '
$ws.Range("D17").Value = '2.045.82'
$ws.Range("E17").Value = '  -2.13%  '
$ws.Range("D18").Value = '36.890.99'
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").Value = '''16.85'
$ws.Range("E19").Value = '  +16.06%  '
$ws.Range("D20").Value = '''75.06'
$ws.Range("E20").Value = '  +2.36%  '
$ws.Range("D21").Value = '0.0₃0905'
$ws.Range("E21").Value = '  +5.92%  '
$ws.Range("E22").Value = '  +3.13%  '
$ws.Range("D23").Value = '''237.20'
$ws.Range("E23").Value = '  -1.71%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").Value = '''2.39'
$ws.Range("E25").Value = '  -3.45%  '
$ws.Range("D26").Value = '''2.28'
$ws.Range("E26").Value = '  +11.22%  '
$ws.Range("D27").Value = '''169.08'
$ws.Range("E27").Value = '  -1.82%  '
$ws.Range("D28").Value = '''9.27'
$ws.Range("E28").Value = '  -0.16%  '
$ws.Range("D29").Value = '''20.20'
$ws.Range("E29").Value = '  -3.47%  '
$ws.Range("D30").Value = '''0.126'
$ws.Range("E30").Value = '  +1.00%  '
$ws.Range("D31").Value = '''1.15'
$ws.Range("E31").Value = '  +4.25%  '
$ws.Range("D32").Value = '''4.74'
$ws.Range("E32").Value = '  +4.00%  '
$ws.Range("D33").Value = '''0.0621'
$ws.Range("E33").Value = '  -1.27%  '
$ws.Range("E34").Value = '  +3.42%  '
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("D36").Value = '''0.0871'
$ws.Range("E36").Value = '  -3.59%  '
$ws.Range("D37").Value = '''2.23'
$ws.Range("E37").Value = '  -1.78%  '
$ws.Range("E38").Value = '  -3.52%  '
$ws.Range("D39").Value = '''0.111'
$ws.Range("E39").Value = '  +14.09%  '
$ws.Range("E40").Value = '  +0.33%  '
$ws.Range("D41").Value = '''17.87'
$ws.Range("E41").Value = '  -0.26%  '
$ws.Range("D42").Value = '''0.0224'
$ws.Range("E42").Value = '  -2.19%  '
$ws.Range("D43").Value = '''1.14'
$ws.Range("E43").Value = '  -3.65%  '
$ws.Range("D44").Value = '''96.96'
$ws.Range("E44").Value = '  -2.90%  '
$ws.Range("D45").Value = '''2.84'
$ws.Range("E45").Value = '  +1.47%  '
$ws.Range("D46").Value = '''4.68'
$ws.Range("E46").Value = '  +15.26%  '
$ws.Range("D47").Value = '''2.48'
$ws.Range("E47").Value = '  +5.18%  '
$ws.Range("D48").Value = '1.283.18'
$ws.Range("E48").Value = '  -3.73%  '
$ws.Range("E49").Value = '  -1.88%  '
$ws.Range("D50").Value = '''6.77'
$ws.Range("E50").Value = '  -3.69%  '
$ws.Range("D51").Value = '2.228.23'
$ws.Range("E51").Value = '  -2.07%  '
